$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1, copying the formatting used by the other header cells (A1:C1)
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Tipo"

# Update existing numeric results for B2 and C2
$ws.Range("B2").Value = 0.2805525182256373
$ws.Range("C2").Value = 0.9946027933975616

# Add new "Tipo" value for the row
$ws.Range("D2").Value = "single"
